$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6450
$ws.Range("F3").Value = 2598
$ws.Range("F5").Value = 1302
$ws.Range("F7").Value = 3184
$ws.Range("F9").Value = 172
$ws.Range("F11").Value = 8045
$ws.Range("F12").Value = 415
$ws.Range("F14").Value = 119
$ws.Range("F15").Value = 18
$ws.Range("F16").Value = 285
$ws.Range("F18").Value = 52
$ws.Range("F21").Value = 10007
$ws.Range("F23").Value = 278
$ws.Range("F24").Value = 34
$ws.Range("F25").Value = 136
$ws.Range("F26").Value = 377
$ws.Range("F27").Value = 39
$ws.Range("F35").Value = 27
$ws.Range("F37").Value = 4015
$ws.Range("F38").Value = 247
$ws.Range("F40").Value = 1932
$ws.Range("F41").Value = 1211
$ws.Range("F42").Value = 135
$ws.Range("F44").Value = 186
$ws.Range("F46").Value = 84
$ws.Range("F47").Value = 66
$ws.Range("F49").Value = 45

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 25
$ws.Range("F13").Value = 21
$ws.Range("F14").Value = 5
$ws.Range("F18").Value = 22
$ws.Range("F20").Value = 17

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6450
$ws.Range("F3").Value = 2598
$ws.Range("F6").Value = 1302
$ws.Range("F8").Value = 3184
$ws.Range("F11").Value = 172
$ws.Range("F13").Value = 8045
$ws.Range("F14").Value = 415
$ws.Range("F16").Value = 119
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 285
$ws.Range("F19").Value = 52
$ws.Range("F22").Value = 10007
$ws.Range("F23").Value = 278
$ws.Range("F24").Value = 34
$ws.Range("F25").Value = 136
$ws.Range("F26").Value = 377
$ws.Range("F27").Value = 39
$ws.Range("F28").Value = 21
$ws.Range("F36").Value = 4015
$ws.Range("F37").Value = 247
$ws.Range("F39").Value = 1933
$ws.Range("F40").Value = 22
$ws.Range("F41").Value = 1211
$ws.Range("F42").Value = 135
$ws.Range("F44").Value = 186
$ws.Range("F46").Value = 84
$ws.Range("F47").Value = 66
$ws.Range("F49").Value = 45
